$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing "Long-term Disability Insurance" / "Medical & Dental
# Insurance" rows for JPMorgan Chase & Co. (positive/neutral counts were
# transposed for these two rows).
$ws.Range("D29").Value = 3
$ws.Range("E29").Value = 1

$ws.Range("D30").Value = 2
$ws.Range("E30").Value = 1

# Append newly scraped JPM mention-count rows (company, name, count,
# positive, neutral, negative).
$newRows = @(
    @("JPMorgan Chase & Co.", "401K",           6, 4, 1, 1),
    @("JPMorgan Chase & Co.", "Life Insurance", 4, 4, 0, 0),
    @("JPMorgan Chase & Co.", "ETFs",           2, 2, 0, 0),
    @("JPMorgan Chase & Co.", "Savings",        3, 2, 0, 1)
)

$startRow = 37
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Range("A$r").Value = $data[0]
    $ws.Range("B$r").Value = $data[1]
    $ws.Range("C$r").Value = $data[2]
    $ws.Range("D$r").Value = $data[3]
    $ws.Range("E$r").Value = $data[4]
    $ws.Range("F$r").Value = $data[5]
}

# Reflect the final selection/scroll state from the saved workbook.
$ws.Range("F29").Select()
